$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-19: updated price/volume figures
$ws.Range("D2").Value = "'67.629.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "'2.419.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'554.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "'161.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.509"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("D9").Value = "'0.157"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.90%  "
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("D11").Value = "'0.326"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.05%  "
$ws.Range("D12").Value = "'4.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "'67.533.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("D14").Value = "'0.0000168"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").Value = "'22.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").Value = "'10.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.64%  "
$ws.Range("D17").Value = "'335.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("D18").Value = "'6.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.88%  "
$ws.Range("D19").Value = "'3.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.31%  "

# Row 20/21: Dai and SuiNetwork swapped order (SuiNetwork now ranks above Dai)
$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D20").Value = "'1.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.83%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.01%  "

# Rows 22-49: updated price/volume figures
$ws.Range("D22").Value = "'66.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("D23").Value = "'3.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("D24").Value = "'8.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").Value = "'0.0₃0809"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("D26").Value = "'7.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "'421.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("D29").Value = "'1.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").Value = "'1.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("D31").Value = "'160.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.71%  "
$ws.Range("D32").Value = "'18.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").Value = "'17.70"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "'0.103"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.25%  "
$ws.Range("D36").Value = "'0.293"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.54%  "
$ws.Range("D37").Value = "'4.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.94%  "
$ws.Range("D38").Value = "'1.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").Value = "'1.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.04%  "
$ws.Range("D40").Value = "'2.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("D41").Value = "'3.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("D42").Value = "'128.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.00%  "
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").Value = "'0.476"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").Value = "'0.554"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("D46").Value = "'0.0913"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").Value = "'1.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.42%  "
$ws.Range("D49").Value = "'16.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.65%  "

# Row 50/51: THORChain and BabyDogeCoin swapped order (BabyDogeCoin now ranks above THORChain)
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.0₆0203"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.51%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'4.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.23%  "

